$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------------
# Change 1: paragraph "2. Il y a trop de choses a savoir dans trop peu de
# temps." - remove the explicit w:sz / w:szCs (24) from both the paragraph
# mark rPr and the run rPr, keeping everything else (incl. rsid attributes)
# identical.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(52)
$rng1 = $p1.Range
$xml1 = '<w:p ' + $wns + ' w14:paraId="1F187B64" w14:textId="77777777" w:rsidR="00AB4FE7" w:rsidRPr="005905F9" w:rsidRDefault="00AB4FE7" w:rsidP="00AB4FE7"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr></w:pPr><w:r w:rsidRPr="005905F9"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve">2. Il y a trop de choses à savoir dans trop peu de temps. </w:t></w:r></w:p>'
[void]$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: paragraph "3. Vous pouvez etre tenus de travailler des heures
# supplementaires." - same sz / szCs removal.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(53)
$rng2 = $p2.Range
$xml2 = '<w:p ' + $wns + ' w14:paraId="09425D50" w14:textId="77777777" w:rsidR="00AB4FE7" w:rsidRDefault="00AB4FE7" w:rsidP="00805410"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr></w:pPr><w:r w:rsidRPr="005905F9"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve">3. Vous pouvez être tenus de travailler des heures supplémentaires. </w:t></w:r></w:p>'
[void]$rng2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 3: merge the two runs of the "Il est vrai que ..." paragraph into a
# single run/sentence, and move the page-break + "Mon conseil..." text into
# its own new paragraph that now carries the <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(56)
$p4 = $d.Paragraphs.Item(57)
$rngStart = $p3.Range.Start
$rngEnd = $p4.Range.End
$rng3 = $d.Range($rngStart, $rngEnd)

$nbsp = [char]0x00A0

$paraA = '<w:p ' + $wns + ' w14:paraId="26CEBE01" w14:textId="77777777" w:rsidR="001408ED" w:rsidRDefault="001408ED" w:rsidP="00657730"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Il est vrai que nos universités n’ayant pas dans leur cursus les dernières mises à jour concernant l’évolution des métiers de développeurs, ne nous préparent absolument pas du tout pour le marché de l’emploi, et je comprends facilement que des étudiants en finissant se sentent un peu déroutés. </w:t></w:r></w:p>'
$paraB = '<w:p ' + $wns + '><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Mon conseil pour les développeurs informatiques e</w:t></w:r><w:r w:rsidR="004D4F7D"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>st' + $nbsp + ':</w:t></w:r></w:p>'

[void]$rng3.InsertXML($paraA + $paraB)

Write-Host "All changes applied."
